$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block in rows 4-8 (columns D, L, M, N, O, P, Q, R, S) is being
# cyclically rotated: old row 6 -> row 4, old row 7 -> row 5, old row 8 -> row 6,
# old row 4 -> row 7, old row 5 -> row 8.

$ws.Range("D4").Value = 44285
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R4").Value = "Provincia del Elquí"
$ws.Range("S4").Value = 1200

$ws.Range("D5").Value = 44285
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 90
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R5").Value = "Provincia del Elquí"
$ws.Range("S5").Value = 1000

$ws.Range("D6").Value = 44285
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 75
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R6").Value = "Provincia del Elquí"
$ws.Range("S6").Value = 800

$ws.Range("D7").Value = 44309
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = "`$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1200

$ws.Range("D8").Value = 44309
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = "`$/caja 15 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1000
